# Apply the ProjectConfiguration.xlsx edit described by the commit:
# "test project from `esqlabsR::initProject()`"
#
# Content-level change: the "dataFile" row's Value cell (B12) changes from
# the old GHM example data file name to the new TestProject example file
# name used by the esqlabsR test-project template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = "TestProject_TimeValuesData.xlsx"

# The saved selection/active cell also moved (cosmetic, matches the diff's
# <selection activeCell="A8" sqref="A8"/>).
$ws.Range("A8").Select()
